$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47 (shifts existing rows 47-143 down to 48-144)
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with its data. Columns A,B,C,E,F,G,H,I,N,O,Q,R are
# unchanged from the surrounding rows (same market/category block), so copy
# them from row 48 (the row that used to be row 47 before the insert). Use
# Value2 for reads (Value round-trips oddly through the COM Variant shim).
$ws.Range("A47").Value = $ws.Range("A48").Value2
$ws.Range("B47").Value = $ws.Range("B48").Value2
$ws.Range("C47").Value = $ws.Range("C48").Value2
$ws.Range("D47").Value = 45002
$ws.Range("E47").Value = $ws.Range("E48").Value2
$ws.Range("F47").Value = $ws.Range("F48").Value2
$ws.Range("G47").Value = $ws.Range("G48").Value2
$ws.Range("H47").Value = $ws.Range("H48").Value2
$ws.Range("I47").Value = $ws.Range("I48").Value2
$ws.Range("J47").Value = 170
$ws.Range("K47").Value = 8000
$ws.Range("L47").Value = 8500
$ws.Range("M47").Value = 8265
$ws.Range("N47").Value = $ws.Range("N48").Value2
$ws.Range("O47").Value = $ws.Range("O48").Value2
$ws.Range("P47").Value = 138
$ws.Range("Q47").Value = $ws.Range("Q48").Value2
$ws.Range("R47").Value = $ws.Range("R48").Value2
